# Regenerate the "K" (strikeouts) column (column G) values for the
# snell_blake 2022 save_data sheet, replacing the old values with the
# newly computed ones.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(6,6,5,6,6,13,7,5,10,8,4,7,10,8,9,7,5,5,11,12,4,5,7,4,6,7,5,4,0)

for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
